$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.284.86"
$ws.Range("E2").Value = "  +4.01%  "

$ws.Range("D3").Value = "3.627.23"
$ws.Range("E3").Value = "  +8.80%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "243.19"
$ws.Range("E5").Value = "  +5.15%  "

$ws.Range("D6").Value = "654.83"
$ws.Range("E6").Value = "  +6.56%  "

$ws.Range("E7").Value = "  +6.01%  "

$ws.Range("D8").Value = "0.406"
$ws.Range("E8").Value = "  +4.52%  "

$ws.Range("D10").Value = "1.02"
$ws.Range("E10").Value = "  +7.97%  "

$ws.Range("D11").Value = "3.628.40"
$ws.Range("E11").Value = "  +8.81%  "

$ws.Range("D12").Value = "44.00"
$ws.Range("E12").Value = "  +3.07%  "

$ws.Range("D13").Value = "0.202"
$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("D14").Value = "6.41"
$ws.Range("E14").Value = "  +3.59%  "

$ws.Range("D15").Value = "4.320.48"
$ws.Range("E15").Value = "  +8.91%  "

$ws.Range("D16").Value = "96.139.71"
$ws.Range("E16").Value = "  +4.04%  "

$ws.Range("D17").Value = "0.0000259"
$ws.Range("E17").Value = "  +6.04%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.626.88"
$ws.Range("E18").Value = "  +8.57%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "12.73"
$ws.Range("E20").Value = "  +14.63%  "

$ws.Range("D21").Value = "18.41"
$ws.Range("E21").Value = "  +6.70%  "

$ws.Range("B22").Value = "Stellar"
$ws.Range("C22").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D22").Value = "0.491"
$ws.Range("E22").Value = "  +13.86%  "

$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "518.12"
$ws.Range("E24").Value = "  +5.35%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000199"
$ws.Range("E25").Value = "  +9.53%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "6.84"
$ws.Range("E26").Value = "  +4.97%  "

$ws.Range("D27").Value = "97.79"
$ws.Range("E27").Value = "  +5.45%  "

$ws.Range("D28").Value = "12.80"
$ws.Range("E28").Value = "  +7.27%  "

$ws.Range("D29").Value = "3.21"
$ws.Range("E29").Value = "  +22.04%  "

$ws.Range("D30").Value = "11.59"
$ws.Range("E30").Value = "  +4.70%  "

$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "0.995"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").Value = "0.178"
$ws.Range("E34").Value = "  +4.21%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "32.09"
$ws.Range("E35").Value = "  +13.80%  "

$ws.Range("D36").Value = "0.568"
$ws.Range("E36").Value = "  +8.67%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "8.06"
$ws.Range("E37").Value = "  +8.18%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "563.55"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("D39").Value = "1.48"
$ws.Range("E39").Value = "  +7.97%  "

$ws.Range("D40").Value = "0.946"
$ws.Range("E40").Value = "  +9.74%  "

$ws.Range("D41").Value = "0.152"
$ws.Range("E41").Value = "  +3.04%  "

$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").Value = "1.77"
$ws.Range("E43").Value = "  +6.67%  "

$ws.Range("D44").Value = "5.86"
$ws.Range("E44").Value = "  +9.13%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "2.31"
$ws.Range("E45").Value = "  +10.53%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0427"
$ws.Range("E46").Value = "  +4.62%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "23.78"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "33.75"
$ws.Range("E48").Value = "  +51.17%  "

$ws.Range("D49").Value = "54.75"
$ws.Range("E49").Value = "  +4.96%  "

$ws.Range("E50").Value = "  +5.34%  "

$ws.Range("D51").Value = "3.49"
$ws.Range("E51").Value = "  -3.34%  "
